# SHKMail made more universal (issue #704) + property support matrix update
#
# Adds a new "Mail" sharer column (L) to the support matrix worksheet and
# marks the relevant property cells supported for it, plus a couple of
# incidental "x" fixes in the Facebook (J) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Mail" column header (row 14) -------------------------------------
$ws.Range("L14").Value = "Mail"

# --- New shared "x" marks in existing Facebook column (J) ------------------
$ws.Range("J20").Value = "x"
$ws.Range("J28").Value = "x"
$ws.Range("J33").Value = "x"

# --- New "Mail" column (L) support marks ------------------------------------
$ws.Range("L16").Value = "x"
$ws.Range("L18").Value = "x"
$ws.Range("L19").Value = "x"
$ws.Range("L20").Value = "x"
$ws.Range("L22").Value = "x"
$ws.Range("L28").Value = "x"
$ws.Range("L30").Value = "x"
$ws.Range("L31").Value = "x"
$ws.Range("L32").Value = "x"
$ws.Range("L34").Value = "x"
$ws.Range("L40").Value = "x"
$ws.Range("L42").Value = "x"
$ws.Range("L43").Value = "x"
$ws.Range("L44").Value = "x"
$ws.Range("L46").Value = "x"
$ws.Range("L52").Value = "x"
$ws.Range("L54").Value = "x"
$ws.Range("L55").Value = "x"
$ws.Range("L56").Value = "x"
$ws.Range("L58").Value = "x"
$ws.Range("L63").Value = "N/A"

# --- View state: re-split panes at row 13 / top-left A14 --------------------
$win = $excel.ActiveWindow
$win.SplitRow = 13
$win.SplitColumn = 0
$win.FreezePanes = $false

$ws.Range("A6").Select() | Out-Null
$ws.Range("L35").Select() | Out-Null

Write-Output "done"
